# Refined metadata to be additional tab
$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# ---------------------------------------------------------------------------
# 1) Update the "panel_query_time" column (F) on the "data" sheet - the query
#    timestamps were refreshed to a later run.
# ---------------------------------------------------------------------------
$newTimes = @(
    "2021-10-05 14:35:31.067704",
    "2021-10-05 14:35:31.067712",
    "2021-10-05 14:35:31.067715",
    "2021-10-05 14:35:31.067718",
    "2021-10-05 14:35:31.067721",
    "2021-10-05 14:35:31.067723",
    "2021-10-05 14:35:31.067726",
    "2021-10-05 14:35:31.067728",
    "2021-10-05 14:35:31.067731",
    "2021-10-05 14:35:31.067734",
    "2021-10-05 14:35:31.067736",
    "2021-10-05 14:35:31.067739",
    "2021-10-05 14:35:31.067741",
    "2021-10-05 14:35:31.067744",
    "2021-10-05 14:35:31.067746",
    "2021-10-05 14:35:31.067749",
    "2021-10-05 14:35:31.067752",
    "2021-10-05 14:35:31.067755",
    "2021-10-05 14:35:31.067757",
    "2021-10-05 14:35:31.067760",
    "2021-10-05 14:35:31.067762",
    "2021-10-05 14:35:31.067765",
    "2021-10-05 14:35:31.067767",
    "2021-10-05 14:35:31.067770"
)

for ($i = 0; $i -lt $newTimes.Length; $i++) {
    $row = $i + 2
    $dataSheet.Cells.Item($row, 6).Value = $newTimes[$i]
}

# ---------------------------------------------------------------------------
# 2) Add a new "metadata" worksheet right after "data", describing the panel
#    query that produced this workbook.
# ---------------------------------------------------------------------------
$metaSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $dataSheet)
$metaSheet.Name = "metadata"

# Header row (bold, bordered, centered - matching the "data" sheet headers)
$metaSheet.Range("B1").Value = "data_name"
$metaSheet.Range("C1").Value = "data_id"
$metaSheet.Range("D1").Value = "data_version"
$metaSheet.Range("E1").Value = "data_version_created"
$metaSheet.Range("F1").Value = "panel_query_time"
$metaSheet.Range("G1").Value = "panel_get_request"

$headerRange = $metaSheet.Range("B1:G1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Row index cell (A2) - same bold/bordered/centered style as "data" sheet's A column
$idxCell = $metaSheet.Range("A2")
$idxCell.Value = 0
$idxCell.Font.Bold = $true
$idxCell.HorizontalAlignment = -4108
$idxCell.VerticalAlignment = -4160
$idxCell.Borders.LineStyle = 1

# Data row
$metaSheet.Range("B2").Value = "Renal abnormalities of calcium and phosphate metabolism"
$metaSheet.Range("C2").Value = 192

$versionCell = $metaSheet.Range("D2")
$versionCell.NumberFormat = "@"
$versionCell.Value = "0.34"

$metaSheet.Range("E2").Value = "2021-08-12T04:23:33.402048Z"
$metaSheet.Range("F2").Value = "2021-10-05 14:35:31.064511"
$metaSheet.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/192/?format=json"

# Leave the focus back on the "data" sheet, as in the original workbook.
$dataSheet.Select()
